# Weekly update: insert the latest week's Perejil price record at the top
# of the data block (row 521), pushing the existing historical rows down
# by one. This mirrors the source system prepending the newest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 521:549 down to 522:550, leaving a blank row 521.
$ws.Rows("521").Insert()

# Populate the new row 521 with this week's record.
$ws.Range("A521").Value2 = 9
$ws.Range("B521").Value2 = 'Vega Central Mapocho de Santiago'
$ws.Range("C521").Value2 = 'Metropolitana'
$ws.Range("D521").Value2 = 45041
$ws.Range("E521").Value2 = 13
$ws.Range("F521").Value2 = 100112044
$ws.Range("G521").Value2 = 'Perejil'
$ws.Range("H521").Value2 = 'Sin especificar'
$ws.Range("I521").Value2 = 'Primera'
$ws.Range("J521").Value2 = 70
$ws.Range("K521").Value2 = 13000
$ws.Range("L521").Value2 = 14000
$ws.Range("M521").Value2 = 13500
$ws.Range("N521").Value2 = '$/docena de atados'
$ws.Range("O521").Value2 = 'Región Metropolitana'
$ws.Range("P521").Value2 = 4500
$ws.Range("Q521").Value2 = 3
$ws.Range("R521").Value2 = 'Hortaliza'
